{"js": "// Each entry is [oldText, newText]; oldText values are unique in the\n// document so a plain text search+replace is unambiguous for every cell.\nconst pairs = [\n  [\"2025-03-31 Monday\", \"2025-04-01 Tuesday\"],\n  [\"967\u00d74=\", \"532\u00d77=\"],\n  [\"794\u00d79=\", \"773\u00d72=\"],\n  [\"674\u00d72=\", \"720\u00d72=\"],\n  [\"206\u00d77=\", \"444\u00d78=\"],\n  [\"543\u00d78=\", \"848\u00d74=\"],\n  [\"452\u00d76=\", \"106\u00d72=\"],\n  [\"508\u00d78=\", \"848\u00d73=\"],\n  [\"171\u00d77=\", \"556\u00d72=\"],\n  [\"453\u00d73=\", \"615\u00d75=\"],\n  [\"725\u00d79=\", \"856\u00d78=\"],\n  [\"977\u00d73=\", \"384\u00d74=\"],\n  [\"367\u00d74=\", \"293\u00d77=\"],\n  [\"598\u00d77=\", \"771\u00d75=\"],\n  [\"623\u00d75=\", \"611\u00d79=\"],\n  [\"610\u00d75=\", \"949\u00d74=\"],\n  [\"262\u00d72=\", \"226\u00d78=\"],\n  [\"619\u00d78=\", \"431\u00d75=\"],\n  [\"122\u00d78=\", \"847\u00d78=\"],\n  [\"982\u00d73=\", \"977\u00d78=\"],\n  [\"350\u00d78=\", \"571\u00d74=\"],\n  [\"678\u00d73=\", \"898\u00d75=\"],\n  [\"756\u00d72=\", \"923\u00d79=\"],\n  [\"333\u00d74=\", \"391\u00d77=\"],\n  [\"562\u00d75=\", \"951\u00d79=\"],\n  [\"349\u00d73=\", \"426\u00d74=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "$d = $word.ActiveDocument\n\n# Each entry is (oldText, newText); oldText values are unique in the\n# document so a plain Find/Replace is unambiguous for every cell.\n$pairs = @(\n  @(\"2025-03-31 Monday\", \"2025-04-01 Tuesday\"),\n  @(\"967\u00d74=\", \"532\u00d77=\"),\n  @(\"794\u00d79=\", \"773\u00d72=\"),\n  @(\"674\u00d72=\", \"720\u00d72=\"),\n  @(\"206\u00d77=\", \"444\u00d78=\"),\n  @(\"543\u00d78=\", \"848\u00d74=\"),\n  @(\"452\u00d76=\", \"106\u00d72=\"),\n  @(\"508\u00d78=\", \"848\u00d73=\"),\n  @(\"171\u00d77=\", \"556\u00d72=\"),\n  @(\"453\u00d73=\", \"615\u00d75=\"),\n  @(\"725\u00d79=\", \"856\u00d78=\"),\n  @(\"977\u00d73=\", \"384\u00d74=\"),\n  @(\"367\u00d74=\", \"293\u00d77=\"),\n  @(\"598\u00d77=\", \"771\u00d75=\"),\n  @(\"623\u00d75=\", \"611\u00d79=\"),\n  @(\"610\u00d75=\", \"949\u00d74=\"),\n  @(\"262\u00d72=\", \"226\u00d78=\"),\n  @(\"619\u00d78=\", \"431\u00d75=\"),\n  @(\"122\u00d78=\", \"847\u00d78=\"),\n  @(\"982\u00d73=\", \"977\u00d78=\"),\n  @(\"350\u00d78=\", \"571\u00d74=\"),\n  @(\"678\u00d73=\", \"898\u00d75=\"),\n  @(\"756\u00d72=\", \"923\u00d79=\"),\n  @(\"333\u00d74=\", \"391\u00d77=\"),\n  @(\"562\u00d75=\", \"951\u00d79=\"),\n  @(\"349\u00d73=\", \"426\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute(\n    [ref]$old,    # FindText\n    [ref]$False,  # MatchCase\n    [ref]$False,  # MatchWholeWord\n    [ref]$False,  # MatchWildcards\n    [ref]$False,  # MatchSoundsLike\n    [ref]$False,  # MatchAllWordForms\n    [ref]$True,   # Forward\n    [ref]1,       # Wrap (wdFindContinue)\n    [ref]$False,  # Format\n    [ref]$new,    # ReplaceWith\n    [ref]2        # Replace (wdReplaceAll)\n  ) | Out-Null\n}\n\n"}
